$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# Delete the empty "Sheet1", keep "Sheet2" which has all the BoM/PnP data.
[void]$wb.Worksheets.Item("Sheet1").Delete()

# Rename the remaining "Sheet2" to "Sheet1".
$wb.Worksheets.Item("Sheet2").Name = "Sheet1"
[void]$wb.Worksheets.Item("Sheet1").Activate()

$excel.DisplayAlerts = $true
